# Insert a new weekly price observation as row 267 in the "Cilantro" sheet.
# This shifts the existing rows 267-296 down to 268-297 (dimension grows
# from A1:R296 to A1:R297) and populates the newly inserted row 267 with
# the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing rows 267..296 down by one to make room for the new record.
$ws.Rows.Item(267).Insert()

# Fill in the new record in row 267.
$ws.Cells.Item(267, 1).Value = 7
$ws.Cells.Item(267, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(267, 3).Value = "Ñuble"
$ws.Cells.Item(267, 4).Value = 45142
$ws.Cells.Item(267, 5).Value = 16
$ws.Cells.Item(267, 6).Value = 100112040
$ws.Cells.Item(267, 7).Value = "Cilantro"
$ws.Cells.Item(267, 8).Value = "Sin especificar"
$ws.Cells.Item(267, 9).Value = "Primera"
$ws.Cells.Item(267, 10).Value = 60
$ws.Cells.Item(267, 11).Value = 1500
$ws.Cells.Item(267, 12).Value = 1500
$ws.Cells.Item(267, 13).Value = 1500
$ws.Cells.Item(267, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(267, 15).Value = "Región de Ñuble"
$ws.Cells.Item(267, 16).Value = 1500
$ws.Cells.Item(267, 17).Value = 1
$ws.Cells.Item(267, 18).Value = "Hortaliza"
